# "Generate Report for Handoff": mark the rows that were waiting on a
# low-priority handoff as handed-off ("ht"), and stamp the new handoff
# timestamp for those rows on each locale sheet.

$wb = $excel.ActiveWorkbook

$zhTimestamp = "2016-08-21 16:40:15"
$deTimestamp = "2016-08-21 16:40:19"

$ws = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $ws.Range("E" + $r).Value = "ht"
    $ws.Range("H" + $r).Value = $zhTimestamp
}

$ws = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $ws.Range("E" + $r).Value = "ht"
    $ws.Range("H" + $r).Value = $deTimestamp
}
